$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing cells per diff ---
$ws.Cells.Item(641, 15).Value = 2   # O641: 0 -> 2
$ws.Cells.Item(643, 18).Value = 0   # R643: blank -> 0
$ws.Cells.Item(644, 18).Value = 0   # R644: blank -> 0

# --- Append new weekly rows 645-653 (columns A-Q; R left blank like source) ---
# Row 645
$ws.Cells.Item(645, 1).Value = 45474
$ws.Cells.Item(645, 2).Value = 3939.949951171875
$ws.Cells.Item(645, 3).Value = 4014.800048828125
$ws.Cells.Item(645, 4).Value = 3809
$ws.Cells.Item(645, 5).Value = 3984.300048828125
$ws.Cells.Item(645, 6).Value = 3984.300048828125
$ws.Cells.Item(645, 7).Value = 1864328
$ws.Cells.Item(645, 8).Value = 2024
$ws.Cells.Item(645, 9).Value = 7
$ws.Cells.Item(645, 10).Value = 1
$ws.Cells.Item(645, 11).Value = 0
$ws.Cells.Item(645, 12).Value = 0
$ws.Cells.Item(645, 13).Value = 0
$ws.Cells.Item(645, 14).Value = 27
$ws.Cells.Item(645, 15).Value = 0
$ws.Cells.Item(645, 16).Value = 0
$ws.Cells.Item(645, 17).Value = 0
$ws.Cells.Item(645, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 646
$ws.Cells.Item(646, 1).Value = 45481
$ws.Cells.Item(646, 2).Value = 3998.89990234375
$ws.Cells.Item(646, 3).Value = 4036.60009765625
$ws.Cells.Item(646, 4).Value = 3680
$ws.Cells.Item(646, 5).Value = 3771.35009765625
$ws.Cells.Item(646, 6).Value = 3771.35009765625
$ws.Cells.Item(646, 7).Value = 3481949
$ws.Cells.Item(646, 8).Value = 2024
$ws.Cells.Item(646, 9).Value = 7
$ws.Cells.Item(646, 10).Value = 8
$ws.Cells.Item(646, 11).Value = 0
$ws.Cells.Item(646, 12).Value = 0
$ws.Cells.Item(646, 13).Value = 0
$ws.Cells.Item(646, 14).Value = 28
$ws.Cells.Item(646, 15).Value = 0
$ws.Cells.Item(646, 16).Value = 0
$ws.Cells.Item(646, 17).Value = 0
$ws.Cells.Item(646, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 647
$ws.Cells.Item(647, 1).Value = 45488
$ws.Cells.Item(647, 2).Value = 3772
$ws.Cells.Item(647, 3).Value = 4093.050048828125
$ws.Cells.Item(647, 4).Value = 3767.75
$ws.Cells.Item(647, 5).Value = 3824.14990234375
$ws.Cells.Item(647, 6).Value = 3824.14990234375
$ws.Cells.Item(647, 7).Value = 2919475
$ws.Cells.Item(647, 8).Value = 2024
$ws.Cells.Item(647, 9).Value = 7
$ws.Cells.Item(647, 10).Value = 15
$ws.Cells.Item(647, 11).Value = 0
$ws.Cells.Item(647, 12).Value = 0
$ws.Cells.Item(647, 13).Value = 0
$ws.Cells.Item(647, 14).Value = 29
$ws.Cells.Item(647, 15).Value = 0
$ws.Cells.Item(647, 16).Value = 0
$ws.Cells.Item(647, 17).Value = 1
$ws.Cells.Item(647, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 648
$ws.Cells.Item(648, 1).Value = 45495
$ws.Cells.Item(648, 2).Value = 3815.39990234375
$ws.Cells.Item(648, 3).Value = 4244
$ws.Cells.Item(648, 4).Value = 3620.300048828125
$ws.Cells.Item(648, 5).Value = 4140
$ws.Cells.Item(648, 6).Value = 4140
$ws.Cells.Item(648, 7).Value = 4182276
$ws.Cells.Item(648, 8).Value = 2024
$ws.Cells.Item(648, 9).Value = 7
$ws.Cells.Item(648, 10).Value = 22
$ws.Cells.Item(648, 11).Value = 0
$ws.Cells.Item(648, 12).Value = 0
$ws.Cells.Item(648, 13).Value = 0
$ws.Cells.Item(648, 14).Value = 30
$ws.Cells.Item(648, 15).Value = 2
$ws.Cells.Item(648, 16).Value = 0
$ws.Cells.Item(648, 17).Value = 0
$ws.Cells.Item(648, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 649
$ws.Cells.Item(649, 1).Value = 45502
$ws.Cells.Item(649, 2).Value = 4175
$ws.Cells.Item(649, 3).Value = 4400
$ws.Cells.Item(649, 4).Value = 4125
$ws.Cells.Item(649, 5).Value = 4330.39990234375
$ws.Cells.Item(649, 6).Value = 4330.39990234375
$ws.Cells.Item(649, 7).Value = 4410312
$ws.Cells.Item(649, 8).Value = 2024
$ws.Cells.Item(649, 9).Value = 7
$ws.Cells.Item(649, 10).Value = 29
$ws.Cells.Item(649, 11).Value = 0
$ws.Cells.Item(649, 12).Value = 0
$ws.Cells.Item(649, 13).Value = 0
$ws.Cells.Item(649, 14).Value = 31
$ws.Cells.Item(649, 15).Value = 0
$ws.Cells.Item(649, 16).Value = 0
$ws.Cells.Item(649, 17).Value = 0
$ws.Cells.Item(649, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 650
$ws.Cells.Item(650, 1).Value = 45509
$ws.Cells.Item(650, 2).Value = 4151.5498046875
$ws.Cells.Item(650, 3).Value = 4433.10009765625
$ws.Cells.Item(650, 4).Value = 4130.0498046875
$ws.Cells.Item(650, 5).Value = 4377.9501953125
$ws.Cells.Item(650, 6).Value = 4377.9501953125
$ws.Cells.Item(650, 7).Value = 3420171
$ws.Cells.Item(650, 8).Value = 2024
$ws.Cells.Item(650, 9).Value = 8
$ws.Cells.Item(650, 10).Value = 5
$ws.Cells.Item(650, 11).Value = 0
$ws.Cells.Item(650, 12).Value = 0
$ws.Cells.Item(650, 13).Value = 0
$ws.Cells.Item(650, 14).Value = 32
$ws.Cells.Item(650, 15).Value = 0
$ws.Cells.Item(650, 16).Value = 0
$ws.Cells.Item(650, 17).Value = 0
$ws.Cells.Item(650, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 651
$ws.Cells.Item(651, 1).Value = 45516
$ws.Cells.Item(651, 2).Value = 4390
$ws.Cells.Item(651, 3).Value = 4669
$ws.Cells.Item(651, 4).Value = 4302
$ws.Cells.Item(651, 5).Value = 4636.35009765625
$ws.Cells.Item(651, 6).Value = 4636.35009765625
$ws.Cells.Item(651, 7).Value = 2330058
$ws.Cells.Item(651, 8).Value = 2024
$ws.Cells.Item(651, 9).Value = 8
$ws.Cells.Item(651, 10).Value = 12
$ws.Cells.Item(651, 11).Value = 0
$ws.Cells.Item(651, 12).Value = 0
$ws.Cells.Item(651, 13).Value = 0
$ws.Cells.Item(651, 14).Value = 33
$ws.Cells.Item(651, 15).Value = 0
$ws.Cells.Item(651, 16).Value = 0
$ws.Cells.Item(651, 17).Value = 0
$ws.Cells.Item(651, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 652
$ws.Cells.Item(652, 1).Value = 45523
$ws.Cells.Item(652, 2).Value = 4630
$ws.Cells.Item(652, 3).Value = 4935.39990234375
$ws.Cells.Item(652, 4).Value = 4630
$ws.Cells.Item(652, 5).Value = 4858.2001953125
$ws.Cells.Item(652, 6).Value = 4858.2001953125
$ws.Cells.Item(652, 7).Value = 2758960
$ws.Cells.Item(652, 8).Value = 2024
$ws.Cells.Item(652, 9).Value = 8
$ws.Cells.Item(652, 10).Value = 19
$ws.Cells.Item(652, 11).Value = 0
$ws.Cells.Item(652, 12).Value = 0
$ws.Cells.Item(652, 13).Value = 0
$ws.Cells.Item(652, 14).Value = 34
$ws.Cells.Item(652, 15).Value = 0
$ws.Cells.Item(652, 16).Value = 0
$ws.Cells.Item(652, 17).Value = 0
$ws.Cells.Item(652, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 653
$ws.Cells.Item(653, 1).Value = 45530
$ws.Cells.Item(653, 2).Value = 4890.89990234375
$ws.Cells.Item(653, 3).Value = 5234.60009765625
$ws.Cells.Item(653, 4).Value = 4836.5
$ws.Cells.Item(653, 5).Value = 5182.7998046875
$ws.Cells.Item(653, 6).Value = 5182.7998046875
$ws.Cells.Item(653, 7).Value = 2462679
$ws.Cells.Item(653, 8).Value = 2024
$ws.Cells.Item(653, 9).Value = 8
$ws.Cells.Item(653, 10).Value = 26
$ws.Cells.Item(653, 11).Value = 0
$ws.Cells.Item(653, 12).Value = 0
$ws.Cells.Item(653, 13).Value = 0
$ws.Cells.Item(653, 14).Value = 35
$ws.Cells.Item(653, 15).Value = 0
$ws.Cells.Item(653, 16).Value = 0
$ws.Cells.Item(653, 17).Value = 0
$ws.Cells.Item(653, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
